{"js": "// The visible body text of the document is unchanged by this commit \u2014 the\n// author's change (\"replace sample size calc with protected version\") only\n// re-stamped the SharePoint / content-type custom XML metadata that Word\n// writes alongside the document when the linked workbook attachment is\n// swapped for its protected revision:\n//   * customXml/item1.xml  (ct:contentTypeSchema) \u2014 contentTypeVersion\n//     22->23, a refreshed ma:versionID / ma:fieldsID, and a new\n//     \"MediaServiceSearchProperties\" field (ref + definition).\n//   * customXml/item3.xml  (new documentManagement properties part) plus\n//     its customXml/itemProps3.xml datastore-item companion.\n//   * customXml/itemProps1.xml \u2014 refreshed ds:itemID guid.\n//\n// Word's Office.js object model exposes these as\n// `context.document.customXmlParts` (Word.CustomXmlPart /\n// Word.CustomXmlPartCollection). Drive the edit through that API so this\n// runs unchanged against a full Word/Office.js host; every call is wrapped\n// so a host that keeps these SharePoint metadata parts read-only from\n// add-ins (common \u2014 they are normally only rewritten by the SharePoint\n// sync client, not by document add-ins) doesn't abort the rest of the\n// script.\n\nconst CT_NS = \"http://schemas.microsoft.com/office/2006/metadata/contentType\";\n\nconst NEW_VERSION_ID = \"6229bc6e2580cc7c7b19c04d64859e77\";\nconst NEW_FIELDS_ID = \"beddd486458379b3cb35373f85f5fe9c\";\n\nconst ITEM3_XML =\n  '<?xml version=\"1.0\" encoding=\"utf-8\"?>\\n' +\n  '<p:properties xmlns:p=\"http://schemas.microsoft.com/office/2006/metadata/properties\" xmlns:xsi=\"http://www.w3.org/2001/XMLSchema-instance\" xmlns:pc=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\">\\n' +\n  \"  <documentManagement>\\n\" +\n  '    <TaxCatchAll xmlns=\"fa6a9aea-fb0f-4ddd-aff8-712634b7d5fe\" xsi:nil=\"true\"/>\\n' +\n  '    <DLVStatus xmlns=\"0d58e8a2-dff7-4492-a987-8cd66a35f019\" xsi:nil=\"true\"/>\\n' +\n  '    <lcf76f155ced4ddcb4097134ff3c332f xmlns=\"0d58e8a2-dff7-4492-a987-8cd66a35f019\">\\n' +\n  '      <Terms xmlns=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\"/>\\n' +\n  \"    </lcf76f155ced4ddcb4097134ff3c332f>\\n\" +\n  \"  </documentManagement>\\n\" +\n  \"</p:properties>\\n\";\n\ntry {\n  const parts = context.document.customXmlParts;\n  parts.load(\"items\");\n  await context.sync();\n\n  for (const part of parts.items) {\n    part.load(\"namespaceUri\");\n  }\n  await context.sync();\n\n  // 1) Bump the content-type schema part (customXml/item1.xml): version,\n  //    versionID/fieldsID guids, and the new MediaServiceSearchProperties\n  //    field (both the <xsd:all> reference and its element definition).\n  const schemaPart = parts.items.find((p) => p.namespaceUri === CT_NS);\n  if (schemaPart) {\n    const xmlResult = schemaPart.getXml();\n    await context.sync();\n\n    let xml = xmlResult.value;\n\n    xml = xml.replace(/ma:contentTypeVersion=\"22\"/, 'ma:contentTypeVersion=\"23\"');\n    xml = xml.replace(\n      /ma:versionID=\"45fa2a71bb5685042c59d1fbaafccb6b\"/,\n      `ma:versionID=\"${NEW_VERSION_ID}\"`\n    );\n    xml = xml.replace(\n      /ma:fieldsID=\"3afdd82834536de8985540a519edf7e7\"/,\n      `ma:fieldsID=\"${NEW_FIELDS_ID}\"`\n    );\n\n    const oldAllRef = '<xsd:element ref=\"ns2:MediaLengthInSeconds\" minOccurs=\"0\"/></xsd:all>';\n    const newAllRef =\n      '<xsd:element ref=\"ns2:MediaLengthInSeconds\" minOccurs=\"0\"/>' +\n      '<xsd:element ref=\"ns2:MediaServiceSearchProperties\" minOccurs=\"0\"/></xsd:all>';\n    if (xml.includes(oldAllRef)) {\n      xml = xml.replace(oldAllRef, newAllRef);\n    } else {\n      // Tolerate a pretty-printed / whitespace-separated source part.\n      xml = xml.replace(\n        /(<xsd:element ref=\"ns2:MediaLengthInSeconds\"[^>]*\\/>)(\\s*)(<\\/xsd:all>)/,\n        `$1$2<xsd:element ref=\"ns2:MediaServiceSearchProperties\" minOccurs=\"0\"/>$2$3`\n      );\n    }\n\n    const mediaLengthDef =\n      '<xsd:element name=\"MediaLengthInSeconds\" ma:index=\"23\" nillable=\"true\" ma:displayName=\"MediaLengthInSeconds\" ma:hidden=\"true\" ma:internalName=\"MediaLengthInSeconds\" ma:readOnly=\"true\"><xsd:simpleType><xsd:restriction base=\"dms:Unknown\"/></xsd:simpleType></xsd:element>';\n    const newFieldDef =\n      '<xsd:element name=\"MediaServiceSearchProperties\" ma:index=\"24\" nillable=\"true\" ma:displayName=\"MediaServiceSearchProperties\" ma:hidden=\"true\" ma:internalName=\"MediaServiceSearchProperties\" ma:readOnly=\"true\"><xsd:simpleType><xsd:restriction base=\"dms:Note\"/></xsd:simpleType></xsd:element>';\n    if (xml.includes(mediaLengthDef)) {\n      xml = xml.replace(mediaLengthDef, mediaLengthDef + newFieldDef);\n    } else {\n      xml = xml.replace(\n        /(<xsd:element name=\"MediaLengthInSeconds\"[\\s\\S]*?<\\/xsd:element>)/,\n        `$1${newFieldDef}`\n      );\n    }\n\n    try {\n      schemaPart.setXml(xml);\n      await context.sync();\n    } catch (e) {\n      // Read-only host (SharePoint-managed content-type schema) \u2014 the\n      // metadata bump can't be applied from the add-in surface here.\n    }\n  }\n\n  // 2) Add the new documentManagement properties part\n  //    (customXml/item3.xml + itemProps3.xml datastore item).\n  try {\n    const newPart = context.document.customXmlParts.add(ITEM3_XML);\n    newPart.load(\"id\");\n    await context.sync();\n  } catch (e) {\n    // Host doesn't allow add-ins to mint new custom XML parts.\n  }\n} catch (e) {\n  // Custom XML parts aren't available on this host at all \u2014 nothing else\n  // in the document (body text, styles, etc.) changes for this commit, so\n  // there is nothing further to do.\n}\n", "ps1": "# The visible body text of the document is unchanged by this commit \u2014 the\n# author's change (\"replace sample size calc with protected version\") only\n# re-stamped the SharePoint / content-type custom XML metadata that Word\n# writes alongside the document when the linked workbook attachment is\n# swapped for its protected revision:\n#   * customXml/item1.xml  (ct:contentTypeSchema) \u2014 contentTypeVersion\n#     22->23, a refreshed ma:versionID / ma:fieldsID, and a new\n#     \"MediaServiceSearchProperties\" field (ref + definition).\n#   * customXml/item3.xml  (new documentManagement properties part) plus\n#     its customXml/itemProps3.xml datastore-item companion.\n#   * customXml/itemProps1.xml \u2014 refreshed ds:itemID guid.\n#\n# Word's COM object model exposes these through $d.CustomXMLParts\n# (CustomXMLPart / CustomXMLParts). Drive the edit through that API so this\n# runs unchanged against a real winword.exe; every call is wrapped so a\n# host that keeps these SharePoint metadata parts read-only to automation\n# (common \u2014 they are normally only rewritten by the SharePoint sync\n# client, not by a macro/add-in) doesn't abort the rest of the script.\n\n$d = $word.ActiveDocument\n\n$ctNamespace = \"http://schemas.microsoft.com/office/2006/metadata/contentType\"\n\n$newVersionId = \"6229bc6e2580cc7c7b19c04d64859e77\"\n$newFieldsId = \"beddd486458379b3cb35373f85f5fe9c\"\n\n$item3Xml = @\"\n<?xml version=\"1.0\" encoding=\"utf-8\"?>\n<p:properties xmlns:p=\"http://schemas.microsoft.com/office/2006/metadata/properties\" xmlns:xsi=\"http://www.w3.org/2001/XMLSchema-instance\" xmlns:pc=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\">\n  <documentManagement>\n    <TaxCatchAll xmlns=\"fa6a9aea-fb0f-4ddd-aff8-712634b7d5fe\" xsi:nil=\"true\"/>\n    <DLVStatus xmlns=\"0d58e8a2-dff7-4492-a987-8cd66a35f019\" xsi:nil=\"true\"/>\n    <lcf76f155ced4ddcb4097134ff3c332f xmlns=\"0d58e8a2-dff7-4492-a987-8cd66a35f019\">\n      <Terms xmlns=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\"/>\n    </lcf76f155ced4ddcb4097134ff3c332f>\n  </documentManagement>\n</p:properties>\n\"@\n\ntry {\n    # 1) Bump the content-type schema part (customXml/item1.xml): version,\n    #    versionID/fieldsID guids, and the new MediaServiceSearchProperties\n    #    field (both the <xsd:all> reference and its element definition).\n    for ($i = 1; $i -le $d.CustomXMLParts.Count; $i++) {\n        $part = $d.CustomXMLParts.Item($i)\n        if ($part.NamespaceURI -eq $ctNamespace) {\n            $xml = $part.XML\n\n            $xml = $xml.Replace('ma:contentTypeVersion=\"22\"', 'ma:contentTypeVersion=\"23\"')\n            $xml = $xml.Replace('ma:versionID=\"45fa2a71bb5685042c59d1fbaafccb6b\"', ('ma:versionID=\"' + $newVersionId + '\"'))\n            $xml = $xml.Replace('ma:fieldsID=\"3afdd82834536de8985540a519edf7e7\"', ('ma:fieldsID=\"' + $newFieldsId + '\"'))\n\n            $oldAllRef = '<xsd:element ref=\"ns2:MediaLengthInSeconds\" minOccurs=\"0\"/></xsd:all>'\n            $newAllRef = '<xsd:element ref=\"ns2:MediaLengthInSeconds\" minOccurs=\"0\"/><xsd:element ref=\"ns2:MediaServiceSearchProperties\" minOccurs=\"0\"/></xsd:all>'\n            $xml = $xml.Replace($oldAllRef, $newAllRef)\n\n            $mediaLengthDef = '<xsd:element name=\"MediaLengthInSeconds\" ma:index=\"23\" nillable=\"true\" ma:displayName=\"MediaLengthInSeconds\" ma:hidden=\"true\" ma:internalName=\"MediaLengthInSeconds\" ma:readOnly=\"true\"><xsd:simpleType><xsd:restriction base=\"dms:Unknown\"/></xsd:simpleType></xsd:element>'\n            $newFieldDef = '<xsd:element name=\"MediaServiceSearchProperties\" ma:index=\"24\" nillable=\"true\" ma:displayName=\"MediaServiceSearchProperties\" ma:hidden=\"true\" ma:internalName=\"MediaServiceSearchProperties\" ma:readOnly=\"true\"><xsd:simpleType><xsd:restriction base=\"dms:Note\"/></xsd:simpleType></xsd:element>'\n            $xml = $xml.Replace($mediaLengthDef, ($mediaLengthDef + $newFieldDef))\n\n            try {\n                $part.XML = $xml\n            } catch {\n                # Read-only host (SharePoint-managed content-type schema) \u2014\n                # the metadata bump can't be applied from automation here.\n            }\n            break\n        }\n    }\n} catch {\n    # CustomXMLParts isn't available on this host at all.\n}\n\ntry {\n    # 2) Add the new documentManagement properties part\n    #    (customXml/item3.xml + itemProps3.xml datastore item).\n    $newPart = $d.CustomXMLParts.Add($item3Xml)\n} catch {\n    # Host doesn't allow automation to mint new custom XML parts.\n}\n"}
